$wb = $excel.ActiveWorkbook

# Update the zh-cn sheet: Correspond Handoff/Handback Datetime for row 5
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-26 12:23:23"
$wsZh.Range("G5").Value = "2016-01-26 12:24:07"

# Update the de-de sheet: Correspond Handoff/Handback Datetime for row 5
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-26 12:23:35"
$wsDe.Range("G5").Value = "2016-01-26 12:24:27"
